$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F99").Value = 90
$ws.Range("G99").Value = 8551.799999999999
$ws.Range("F113").Value = 95
$ws.Range("G113").Value = 6680.4
$ws.Range("F114").Value = 50
$ws.Range("G114").Value = 5756
$ws.Range("F116").Value = 82
$ws.Range("G116").Value = 11065.9
$ws.Range("F123").Value = 77
$ws.Range("G123").Value = 3598.98
$ws.Range("F126").Value = 67
$ws.Range("G126").Value = 16486.69
$ws.Range("F133").Value = 46
$ws.Range("G133").Value = 2325.3
$ws.Range("F140").Value = 39
$ws.Range("G140").Value = 1743.3
$ws.Range("B143").Value = 341633.51
$ws.Range("F199").Value = 298
$ws.Range("G199").Value = 5864.64
$ws.Range("F203").Value = 71
$ws.Range("G203").Value = 2273.42
$ws.Range("B205").Value = 34534.36
$ws.Range("B213").Value = 57756
$ws.Range("B214").Value = 53925
$ws.Range("F215").Value = 25
$ws.Range("G215").Value = 1147.5
$ws.Range("B216").Value = 6036.87
$ws.Range("F259").Value = 29
$ws.Range("G259").Value = 2175.58
$ws.Range("B264").Value = 18701.12
$ws.Range("F291").Value = 1
$ws.Range("G291").Value = 34.55
$ws.Range("B295").Value = 895.91
$ws.Range("F299").Value = 25
$ws.Range("G299").Value = 1930.5
$ws.Range("F300").Value = 127
$ws.Range("G300").Value = 8751.57
$ws.Range("B303").Value = 22628.3
$ws.Range("F315").Value = 59
$ws.Range("G315").Value = 8466.5
$ws.Range("B323").Value = 42025.25
$ws.Range("F328").Value = 5
$ws.Range("G328").Value = 3382.6
$ws.Range("B335").Value = 35326.26
$ws.Range("F365").Value = 9
$ws.Range("G365").Value = 1031.04
$ws.Range("F373").Value = 4
$ws.Range("G373").Value = 260.24
$ws.Range("F383").Value = 11
$ws.Range("G383").Value = 5139.09
$ws.Range("F384").Value = 46
$ws.Range("G384").Value = 3941.28
$ws.Range("B386").Value = 150798.7
$ws.Range("F394").Value = 36
$ws.Range("G394").Value = 9602.280000000001
$ws.Range("F395").Value = 112
$ws.Range("G395").Value = 11994.08
$ws.Range("F398").Value = 2
$ws.Range("G398").Value = 147.4
$ws.Range("F400").Value = 14
$ws.Range("G400").Value = 2221.38
$ws.Range("F405").Value = 205
$ws.Range("G405").Value = 23417.15
$ws.Range("F409").Value = 173
$ws.Range("G409").Value = 23713.11
$ws.Range("F410").Value = 141
$ws.Range("G410").Value = 20343.48
$ws.Range("F415").Value = 172
$ws.Range("G415").Value = 8827.040000000001
$ws.Range("F434").Value = 53
$ws.Range("G434").Value = 5516.77
$ws.Range("F440").Value = 82
$ws.Range("G440").Value = 16284.38
$ws.Range("F445").Value = 118
$ws.Range("G445").Value = 13998.34
$ws.Range("F446").Value = 259
$ws.Range("G446").Value = 15314.67
$ws.Range("F451").Value = 618
$ws.Range("G451").Value = 13404.42
$ws.Range("F452").Value = 197
$ws.Range("G452").Value = 1183.97
$ws.Range("F453").Value = 9
$ws.Range("G453").Value = 747.9
$ws.Range("F456").Value = 96
$ws.Range("G456").Value = 19177.92
$ws.Range("F460").Value = 151
$ws.Range("G460").Value = 22826.67
$ws.Range("F466").Value = 0
$ws.Range("G466").Value = 0
$ws.Range("B471").Value = 614159.79
$ws.Range("F473").Value = 30
$ws.Range("G473").Value = 4602.3
$ws.Range("B486").Value = 59294.89
$ws.Range("F500").Value = 93
$ws.Range("G500").Value = 1896.27
$ws.Range("B504").Value = 7319.28
$ws.Range("F538").Value = 156
$ws.Range("G538").Value = 4639.44
$ws.Range("F540").Value = 226
$ws.Range("G540").Value = 21831.6
$ws.Range("B546").Value = 73261.77
$ws.Range("F561").Value = 36
$ws.Range("G561").Value = 10570.32
$ws.Range("F563").Value = 114
$ws.Range("G563").Value = 5403.6
$ws.Range("F568").Value = 112
$ws.Range("G568").Value = 3724
$ws.Range("B575").Value = 77879.25999999999
$ws.Range("F593").Value = 695
$ws.Range("G593").Value = 8902.950000000001
$ws.Range("F608").Value = 710
$ws.Range("G608").Value = 10458.3
$ws.Range("B609").Value = 122374.09
$ws.Range("F676").Value = 734
$ws.Range("G676").Value = 14569.9
$ws.Range("B681").Value = 49335.95
$ws.Range("F696").Value = 63
$ws.Range("G696").Value = 3287.34
$ws.Range("B705").Value = 43814.87
$ws.Range("F716").Value = 67
$ws.Range("G716").Value = 6599.5
$ws.Range("B722").Value = 76968.73
$ws.Range("F731").Value = 55
$ws.Range("G731").Value = 2196.7
$ws.Range("B743").Value = 13898.65
$ws.Range("F747").Value = 52
$ws.Range("G747").Value = 5486.52
$ws.Range("F748").Value = 53
$ws.Range("G748").Value = 5505.11
$ws.Range("F755").Value = 278
$ws.Range("G755").Value = 7764.54
$ws.Range("F762").Value = 100
$ws.Range("G762").Value = 12283
$ws.Range("F763").Value = 9
$ws.Range("G763").Value = 1579.23
$ws.Range("F764").Value = 35
$ws.Range("G764").Value = 2905.7
$ws.Range("B765").Value = 89216.19
$ws.Range("F801").Value = 136
$ws.Range("G801").Value = 4502.96
$ws.Range("F807").Value = 86
$ws.Range("G807").Value = 3763.36
$ws.Range("B808").Value = 53082.81
$ws.Range("F837").Value = 51
$ws.Range("G837").Value = 4679.25
$ws.Range("F842").Value = 109
$ws.Range("G842").Value = 8761.42
$ws.Range("B843").Value = 28082.71
$ws.Range("F868").Value = 22
$ws.Range("G868").Value = 1469.38
$ws.Range("B870").Value = 8711.66
$ws.Range("F890").Value = 95
$ws.Range("G890").Value = 10596.3
$ws.Range("F895").Value = 241
$ws.Range("G895").Value = 16763.96
$ws.Range("F902").Value = 97
$ws.Range("G902").Value = 13968
$ws.Range("F903").Value = 245
$ws.Range("G903").Value = 29573.95
$ws.Range("B905").Value = 117956.88
$ws.Range("F908").Value = 13
$ws.Range("G908").Value = 2246.01
$ws.Range("F927").Value = 183
$ws.Range("G927").Value = 18832.53
$ws.Range("F932").Value = 20
$ws.Range("G932").Value = 636.2
$ws.Range("B937").Value = 72831.35000000001
$ws.Range("F940").Value = 131
$ws.Range("G940").Value = 4899.4
$ws.Range("F941").Value = 48
$ws.Range("G941").Value = 922.08
$ws.Range("F944").Value = 55
$ws.Range("G944").Value = 1713.8
$ws.Range("F945").Value = 223
$ws.Range("G945").Value = 8340.200000000001
$ws.Range("B946").Value = 30339.22
$ws.Range("F979").Value = 4
$ws.Range("G979").Value = 3799.76
$ws.Range("F986").Value = 9
$ws.Range("G986").Value = 4037.31
$ws.Range("B988").Value = 155568.64
$ws.Range("F1005").Value = 280
$ws.Range("G1005").Value = 21599.2
$ws.Range("F1008").Value = 107
$ws.Range("G1008").Value = 13758.06
$ws.Range("B1009").Value = 559988.46
$ws.Range("B1016").Value = 4378946.14
$ws.Range("B1017").Value = 4378946.14
